$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "34.384.32"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.81%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.789.15"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.43%  "

$ws.Range("E4").Value = "  -0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "226.24"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.39%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.555"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.69%  "

$ws.Range("E7").Value = "  -0.07%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "32.82"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.04%  "

$ws.Range("E9").Value = "  +0.87%  "

$ws.Range("E10").Value = "  +0.40%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0947"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.03%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.047.62"
$ws.Range("D12").Style = "Normal"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.13"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.03%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.785.60"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.65%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.634"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.83%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "34.355.19"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.77%  "

$ws.Range("E17").Value = "  +2.77%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.53"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.39%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "245.26"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.26%  "

$ws.Range("E20").Value = "  +0.83%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.21"
$ws.Range("D21").Style = "Normal"

$ws.Range("E23").Value = "  +1.35%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "167.70"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.56%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.06"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.31%  "

$ws.Range("E26").Value = "  +2.88%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.56"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.87%  "

$ws.Range("E28").Value = "  +1.43%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.01"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +7.81%  "

$ws.Range("E31").Value = "  +1.84%  "

$ws.Range("E32").Value = "  +2.31%  "

$ws.Range("E33").Value = "  +0.49%  "

$ws.Range("E34").Value = "  +1.76%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.61"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.22%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.411.12"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.70%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.685"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.16%  "

$ws.Range("E38").Value = "  +3.00%  "

$ws.Range("E39").Value = "  +0.24%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "84.17"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.49%  "

$ws.Range("E41").Value = "  +2.70%  "

$ws.Range("E42").Value = "  +0.62%  "

$ws.Range("E43").Value = "  +2.66%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.85"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.99%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0528"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.53%  "

$ws.Range("E46").Value = "  +3.13%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.06"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.16%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.948.12"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.50%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "105.40"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.70%  "

$ws.Range("E50").Value = "  -0.08%  "

$ws.Range("E51").Value = "  -2.77%  "
